# Updated the email services
# New form submissions (Newsletter Subscriptions, Contact Forms, Demo Requests,
# Talk to Sales) have come in since the last export; insert them at the top of
# each sheet's data (directly under the header row), pushing existing rows
# down, and widen a few columns whose content got longer.

$wb = $excel.ActiveWorkbook

# Excel stores column widths in "Maximum Digit Width" character units that are
# offset by ~5/6 of a character from the simple ColumnWidth number you assign
# via COM. Subtracting that constant before assignment makes the persisted
# <col width="..."> come out as the clean integer we actually want.
$widthFudge = 0.8333333333333

function Set-ColWidth {
    param($ws, $colIndex, $target)
    $ws.Columns.Item($colIndex).ColumnWidth = $target - $widthFudge
}

function Insert-Rows {
    param($ws, $firstRow, $lastRow)
    $ws.Range("${firstRow}:${lastRow}").Insert()
    $ws.Range("A${firstRow}:Z${lastRow}").ClearFormats()
}

# The exported source data is plain text (phone numbers, ISO dates, ids...)
# even when it looks numeric — leading zeros, a leading "+", etc. all need to
# survive. Plain `.Value = "0123"` lets Excel's usual General-format
# auto-detection reinterpret it as a number (and a bare "YYYY-MM-DD" string
# as a date), so force those through as literal text the same way typing a
# leading apostrophe into a cell does.
function Set-Text {
    param($ws, $addr, $val)
    if ($val -eq "") {
        $ws.Range($addr).Value = ""
    } else {
        $ws.Range($addr).Value = "'" + $val
    }
}

# ---------------------------------------------------------------------------
# Sheet: Newsletter Subscriptions  (A1:C11 -> A1:C13, +2 new rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Newsletter Subscriptions")

Insert-Rows $ws 2 3

$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "ra147001y@gmail.com"
$ws.Range("C2").Value = "2026-01-26 21:26:06"

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "ahmedyaqoobbusiness@gmail.com"
$ws.Range("C3").Value = "2026-01-26 01:45:04"

# ---------------------------------------------------------------------------
# Sheet: Contact Forms  (A1:H4 -> A1:H7, +3 new rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contact Forms")

Insert-Rows $ws 2 4

$ws.Range("A2").Value = 15
$ws.Range("B2").Value = "Michael"
$ws.Range("C2").Value = "Mandapati"
$ws.Range("D2").Value = "michael@warpandweft.com"
Set-Text $ws "E2" "2125460944"
$ws.Range("F2").Value = "Warp & Weft"
$ws.Range("G2").Value = "Hello:  I have been recommended by my CPA, Farid Padela and I'm interested in knowing more about your system. Please call to discuss further.`nThank you,`nMichael`n"
$ws.Range("H2").Value = "2026-01-28 10:37:14"

$ws.Range("A3").Value = 14
$ws.Range("B3").Value = "Rizwan"
$ws.Range("C3").Value = "Mahmood"
$ws.Range("D3").Value = "rizwan@vcs.com.pk"
Set-Text $ws "E3" "03318406191"
$ws.Range("F3").Value = "VCS"
$ws.Range("G3").Value = "Testing of the Message field on General Inquiry Form"
$ws.Range("H3").Value = "2026-01-26 02:31:40"

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "Ahmed"
$ws.Range("C4").Value = "Yaqoob"
$ws.Range("D4").Value = "ahmedyaqoobbusiness@gmail.com"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "test"
$ws.Range("H4").Value = "2026-01-26 01:23:06"

Set-ColWidth $ws 4 31
Set-ColWidth $ws 7 50

# ---------------------------------------------------------------------------
# Sheet: Demo Requests  (A1:I3 -> A1:I10, +7 new rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Demo Requests")

Insert-Rows $ws 2 8

$ws.Range("A2").Value = 13
$ws.Range("B2").Value = "Rizwan"
$ws.Range("C2").Value = "Mahmood"
$ws.Range("D2").Value = "rizwan@vcs.com.pk"
Set-Text $ws "E2" "03318406191"
$ws.Range("F2").Value = "VCS"
Set-Text $ws "G2" "2026-01-27"
$ws.Range("H2").Value = "Testing of Additional Information on Demo Request Form of SPARS Website."
$ws.Range("I2").Value = "2026-01-26 02:28:33"

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "ahmed"
$ws.Range("C3").Value = "yaqoob"
$ws.Range("D3").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "E3" "03147261655"
$ws.Range("F3").Value = "asdf"
Set-Text $ws "G3" "2026-01-29"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = "2026-01-26 01:20:34"

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "ahmed"
$ws.Range("C4").Value = "yaqoob"
$ws.Range("D4").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "E4" "0314247610"
$ws.Range("F4").Value = "vcs"
Set-Text $ws "G4" "2026-01-28"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "2026-01-26 01:15:32"

$ws.Range("A5").Value = 9
$ws.Range("B5").Value = "ahmed"
$ws.Range("C5").Value = "yaqoob"
$ws.Range("D5").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "E5" "03147261655"
$ws.Range("F5").Value = "vcs"
Set-Text $ws "G5" "2026-01-27"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "2026-01-26 01:08:07"

$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "ahmed"
$ws.Range("C6").Value = "yaqoob"
$ws.Range("D6").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "E6" "03147261655"
$ws.Range("F6").Value = "vcs"
Set-Text $ws "G6" "2026-01-29"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = "2026-01-26 00:00:23"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "ahmed"
$ws.Range("C7").Value = "yaqoob"
$ws.Range("D7").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "E7" "03147261655"
$ws.Range("F7").Value = "vcsa"
Set-Text $ws "G7" "2026-01-30"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = "2026-01-25 23:57:57"

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "ahmed"
$ws.Range("C8").Value = "yaqoob"
$ws.Range("D8").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "E8" "03147261655"
$ws.Range("F8").Value = "VCS"
Set-Text $ws "G8" "2026-01-27"
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = "2026-01-25 23:47:43"

Set-ColWidth $ws 4 31
Set-ColWidth $ws 8 50

# ---------------------------------------------------------------------------
# Sheet: Talk to Sales  (A1:L4 -> A1:L6, +2 new rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Talk to Sales")

Insert-Rows $ws 2 3

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Rizwan Mahmood"
$ws.Range("C2").Value = "rizwan@vcs.com.pk"
Set-Text $ws "D2" "+923318406191"
$ws.Range("E2").Value = "VCS"
$ws.Range("F2").Value = "Testing of Message Field on Talk to Sales Form"
$ws.Range("G2").Value = "Testing of Current ERP System Field"
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = "Testing of Specific Requirements or Challenges field"
$ws.Range("K2").Value = "testing of Implementation Timeline field"
$ws.Range("L2").Value = "2026-01-26 02:29:45"

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "ahmed"
$ws.Range("C3").Value = "ahmedyaqoobbusiness@gmail.com"
Set-Text $ws "D3" "03147260655"
$ws.Range("E3").Value = "vcs"
$ws.Range("F3").Value = "test"
$ws.Range("G3").Value = "SAP"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "6 months"
$ws.Range("L3").Value = "2026-01-26 01:21:32"

Set-ColWidth $ws 2 16
Set-ColWidth $ws 4 15
Set-ColWidth $ws 6 48
Set-ColWidth $ws 7 37
Set-ColWidth $ws 10 50
Set-ColWidth $ws 11 42

Write-Host "edit complete"
